$d = $word.ActiveDocument

# Update the date heading
$d.Paragraphs.Item(1).Range.Find.Execute(
    "2025-12-09 Tuesday", $false, $false, $false, $false, $false,
    $true, 1, $false, "2025-12-10 Wednesday", 2) | Out-Null

# Update the division expressions in the table.
# Each cell's Find is scoped to that cell's own Range and uses
# wdReplaceOne (1) so that duplicate expressions elsewhere in the
# table are left untouched.
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "89÷3="; New = "57÷3=" },
    @{ Row = 1;  Col = 2; Old = "96÷2="; New = "24÷9=" },
    @{ Row = 1;  Col = 3; Old = "32÷8="; New = "47÷9=" },
    @{ Row = 1;  Col = 4; Old = "18÷7="; New = "26÷9=" },
    @{ Row = 1;  Col = 5; Old = "49÷5="; New = "33÷4=" },

    @{ Row = 5;  Col = 1; Old = "43÷3="; New = "15÷6=" },
    @{ Row = 5;  Col = 2; Old = "84÷6="; New = "88÷8=" },
    @{ Row = 5;  Col = 3; Old = "69÷4="; New = "12÷9=" },
    @{ Row = 5;  Col = 4; Old = "57÷2="; New = "54÷2=" },
    @{ Row = 5;  Col = 5; Old = "98÷3="; New = "54÷8=" },

    @{ Row = 9;  Col = 1; Old = "88÷3="; New = "77÷5=" },
    @{ Row = 9;  Col = 2; Old = "27÷6="; New = "62÷2=" },
    @{ Row = 9;  Col = 3; Old = "87÷4="; New = "78÷5=" },
    @{ Row = 9;  Col = 4; Old = "98÷3="; New = "64÷5=" },
    @{ Row = 9;  Col = 5; Old = "94÷2="; New = "60÷7=" },

    @{ Row = 13; Col = 1; Old = "89÷3="; New = "91÷5=" },
    @{ Row = 13; Col = 2; Old = "76÷4="; New = "85÷8=" },
    @{ Row = 13; Col = 3; Old = "48÷3="; New = "12÷7=" },
    @{ Row = 13; Col = 4; Old = "44÷6="; New = "44÷2=" },
    @{ Row = 13; Col = 5; Old = "10÷9="; New = "76÷3=" },

    @{ Row = 17; Col = 1; Old = "39÷6="; New = "79÷2=" },
    @{ Row = 17; Col = 2; Old = "98÷9="; New = "30÷5=" },
    @{ Row = 17; Col = 3; Old = "46÷6="; New = "38÷5=" },
    @{ Row = 17; Col = 4; Old = "51÷7="; New = "83÷9=" },
    @{ Row = 17; Col = 5; Old = "45÷7="; New = "43÷4=" }
)

foreach ($r in $replacements) {
    $cellRange = $t.Cell($r.Row, $r.Col).Range
    $cellRange.Find.Execute(
        $r.Old, $false, $false, $false, $false, $false,
        $true, 1, $false, $r.New, 1) | Out-Null
}
